$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6: update Operator (Airtel -> Idea), Country Code (30 -> 91),
# and Mobile Number (7888777888 -> 7405383061).
for ($r = 2; $r -le 6; $r++) {

    # Operator column (I) - brand entry text change
    $ws.Cells.Item($r, 9).Value = "Idea"

    # Country Code column (D) - keep it a genuine number, not text,
    # even though the column is formatted as Text (numFmt "@").
    $dCell = $ws.Cells.Item($r, 4)
    $dFmt = $dCell.NumberFormat
    $dCell.NumberFormat = "General"
    $dCell.Value = 91
    $dCell.NumberFormat = $dFmt

    # Mobile Number column (E) - same numeric-preservation trick.
    $eCell = $ws.Cells.Item($r, 5)
    $eFmt = $eCell.NumberFormat
    $eCell.NumberFormat = "General"
    $eCell.Value = 7405383061
    $eCell.NumberFormat = $eFmt
}
